# Auto update Excel log
# Appends new sensor-log rows to the "PIR" and "mmWave" sheets, matching the
# author's commit: a recovery detection row on PIR, and a fall-detection +
# recovery sequence on mmWave.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append row 126 -------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirRow = 126
$pirData = @("2026-01-30", "14:05:11", "14:00", "Living Room", "RECOVERY_DETECTION", "Inactive")

# Column A holds a date-like string ("2026-01-30"); force text formatting so
# Excel stores it as a literal string instead of auto-converting it to a
# date serial number (matches how every other row in the sheet is stored).
$pir.Cells.Item($pirRow, 1).NumberFormat = "@"
$pir.Cells.Item($pirRow, 1).Value = $pirData[0]
$pir.Cells.Item($pirRow, 2).Value = $pirData[1]
$pir.Cells.Item($pirRow, 3).Value = $pirData[2]
$pir.Cells.Item($pirRow, 4).Value = $pirData[3]
$pir.Cells.Item($pirRow, 5).Value = $pirData[4]
$pir.Cells.Item($pirRow, 6).Value = $pirData[5]

# --- mmWave sheet: append rows 47-49 --------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwaveRows = @(
    @("2026-01-30", "14:04:35", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY"),
    @("2026-01-30", "14:05:11", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:05:21", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 47
for ($i = 0; $i -lt $mmwaveRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $mmwaveRows[$i]

    $mmwave.Cells.Item($r, 1).NumberFormat = "@"
    $mmwave.Cells.Item($r, 1).Value = $rowData[0]
    $mmwave.Cells.Item($r, 2).Value = $rowData[1]
    $mmwave.Cells.Item($r, 3).Value = $rowData[2]
    $mmwave.Cells.Item($r, 4).Value = $rowData[3]
    $mmwave.Cells.Item($r, 5).Value = $rowData[4]
    $mmwave.Cells.Item($r, 6).Value = $rowData[5]
}
